$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the data row for id=5 (Мулькин ...), which is worksheet row 4
$ws.Rows.Item(4).Delete()

# After the first deletion, the row that was id=7 (Холостов ...) shifts up
# from worksheet row 6 to worksheet row 5. Delete it too.
$ws.Rows.Item(5).Delete()

# The remaining trailing row (previously row 7, now row 5) holds just a
# timestamp string in column A. Update its value.
$ws.Range("A5").Value = "2019-08-03 00:00:00"
